$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (Excel's ColumnWidth uses a different scale than the
# stored OOXML "width" attribute; subtracting 6/7 compensates for the
# padding offset so the saved width matches the target integer exactly)
$ws.Columns.Item(3).ColumnWidth = 53 - 6/7
$ws.Columns.Item(4).ColumnWidth = 36 - 6/7
$ws.Columns.Item(6).ColumnWidth = 16 - 6/7
$ws.Columns.Item(8).ColumnWidth = 21 - 6/7

# Opportunity IDs in column A are stored as text (not numbers); force text
# formatting on these cells before writing the new numeric-looking IDs so
# they stay text, matching the source data.
$ws.Range("A2:A6").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = "1329367"
$ws.Range("B2").Value = "https://aiesec.org/opportunity/global-talent/1329367"
$ws.Range("C2").Value = "Sales Development Intern"
$ws.Range("D2").Value = "Navi Mumbai, Maharashtra, India"
$ws.Range("F2").Value = "0 applicants"
$ws.Range("G2").Value = "6 - 18 Months"
$ws.Range("H2").Value = "AHY CONSULTING LLP"

# Row 3
$ws.Range("A3").Value = "1328685"
$ws.Range("B3").Value = "https://aiesec.org/opportunity/global-talent/1328685"
$ws.Range("C3").Value = "Medical Advisor (Russian Speaker)"
$ws.Range("D3").Value = "İstanbul, Türkiye"
$ws.Range("F3").Value = "5 applicants"
$ws.Range("H3").Value = "International Plus"

# Row 4
$ws.Range("A4").Value = "1325297"
$ws.Range("B4").Value = "https://aiesec.org/opportunity/global-talent/1325297"
$ws.Range("C4").Value = "International Sales Representetive Spanish Speaker"
$ws.Range("D4").Value = "Maslak, Sarıyer/İstanbul, Türkiye"
$ws.Range("F4").Value = "24 applicants"
$ws.Range("H4").Value = "Esvita Clinic"

# Row 5
$ws.Range("A5").Value = "1321054"
$ws.Range("B5").Value = "https://aiesec.org/opportunity/global-talent/1321054"
$ws.Range("C5").Value = "International Sales Representetive Russian Speaker"
$ws.Range("D5").Value = "Maslak, Sarıyer/İstanbul, Türkiye"
$ws.Range("F5").Value = "15 applicants"
$ws.Range("G5").Value = "6 - 18 Months"
$ws.Range("H5").Value = "Esvita Clinic"

# Row 6
$ws.Range("A6").Value = "1321053"
$ws.Range("B6").Value = "https://aiesec.org/opportunity/global-talent/1321053"
$ws.Range("C6").Value = "International Sales Representetive German Speaker"
$ws.Range("D6").Value = "Maslak, Sarıyer/İstanbul, Türkiye"
$ws.Range("F6").Value = "17 applicants"
$ws.Range("G6").Value = "6 - 18 Months"
$ws.Range("H6").Value = "Esvita Clinic"
